# bouldering pushed to production + update hangboard progress
#
# Appends one new hangboard-session row (row 41) to Sheet1, mirroring the
# layout/style of the existing data rows, and moves the active selection
# to the new last cell (J41) just like the author's workbook showed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 41 so it inherits the number format / style of the
# row above it (row 40), matching the workbook's existing "General"-style
# data rows instead of picking up the blank default style.
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "17 Mar 2021"
$ws.Range("B41").Value = "0,6"
$ws.Range("C41").Value = "-7.5,6"
$ws.Range("D41").Value = "-17.5,6"
$ws.Range("E41").Value = "-37.5,4,2,6"
$ws.Range("F41").Value = "-27.5,6"
$ws.Range("G41").Value = "-12.5,6"
$ws.Range("H41").Value = "-37.5,4,8,7"
$ws.Range("I41").Value = "-20,5,9"
$ws.Range("J41").Value = "-27.5,5,8"

$ws.Range("J41").Select() | Out-Null
